# Sample Project / Main.xlsx — row 11 (R40 rule) "From" value is re-entered
# as the text "1" (kept as text, not a number) while leaving every other
# cell/style untouched.
#
# A plain $ws.Range("B11").Value = "1" gets auto-coerced to a Number by
# Excel's usual "looks like a number" inference, and forcing text via
# NumberFormat="@" (or a leading apostrophe) bakes a new quote-prefixed
# cell style into the workbook. Routing the literal through a formula and
# then collapsing it back to a value with PasteSpecial (xlPasteValues)
# gives a plain text cell without disturbing the existing style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.Formula = "=CHAR(49)"
$cell.Copy()
$cell.PasteSpecial(-4163)
